# Carry marks update: fill in the previously-missing "TEST (30%)" scores
# (column E) for a handful of students. The "TOTAL (50%)" column (G) holds
# a shared formula (=D/30*10+E/50*30+F/50*20) so it recalculates on its own
# once E is populated.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value  = 11
$ws.Range("E11").Value = 23.5
$ws.Range("E21").Value = 23
$ws.Range("E23").Value = 26

# Make sure the dependent totals in column G are up to date.
$excel.Calculate()

# Match the author's final view state: scrolled down a bit further and with
# cell E22 (one of the newly-edited rows) active.
$ws.Activate()
$ws.Range("E22").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
